# ============================================================
# Add 2022-Q4 data:
#  1. Insert a new worksheet "2022-Q4" right after "总计" (summary sheet),
#     holding per-fund holdings data for the new quarter.
#  2. Insert a new summary row in "总计" for 2022-Q4 (count=14, value=1.72),
#     shifting the existing quarterly rows down by one.
# ============================================================

$wb = $excel.ActiveWorkbook

# ---- locate the "总计" (summary) sheet; it is the first sheet ----
$summary = $wb.Worksheets.Item(1)

# ---- 1. Create the new "2022-Q4" sheet right after "总计" ----
$newSheet = $wb.Worksheets.Add($null, $summary)
$newSheet.Name = "2022-Q4"

# Use the next sheet (the former "2022-Q3", now pushed to position 3) as a
# formatting template: bold/centered/bordered header row + index column.
$template = $wb.Worksheets.Item(3)
$template.Range("B1:H1").Copy()
$newSheet.Range("B1").PasteSpecial(-4122)
$template.Range("A2:A15").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)

# Fund code / scale / position / ratio / value columns are stored as text
# (to preserve leading zeros and fixed decimal formatting) in this workbook.
$newSheet.Range("B2:B15").NumberFormat = "@"
$newSheet.Range("D2:G15").NumberFormat = "@"

# ---- header row ----
$newSheet.Range('B1').Value = '基金代码'
$newSheet.Range('C1').Value = '基金名称'
$newSheet.Range('D1').Value = '基金规模'
$newSheet.Range('E1').Value = '股票总仓位'
$newSheet.Range('F1').Value = '仓位占比'
$newSheet.Range('G1').Value = '持有市值(亿元)'
$newSheet.Range('H1').Value = '仓位排名'

# ---- fund holdings detail rows ----
$newSheet.Range('A2').Value = 0
$newSheet.Range('B2').Value = '005354'
$newSheet.Range('C2').Value = '富国沪港深行业精选灵活配置混合A'
$newSheet.Range('D2').Value = '30.68'
$newSheet.Range('E2').Value = '79.74'
$newSheet.Range('F2').Value = '2.37'
$newSheet.Range('G2').Value = '0.7271'
$newSheet.Range('H2').Value = 8
$newSheet.Range('A3').Value = 1
$newSheet.Range('B3').Value = '012434'
$newSheet.Range('C3').Value = '银华多元回报一年持有期混合'
$newSheet.Range('D3').Value = '20.31'
$newSheet.Range('E3').Value = '87.60'
$newSheet.Range('F3').Value = '1.82'
$newSheet.Range('G3').Value = '0.3696'
$newSheet.Range('H3').Value = 5
$newSheet.Range('A4').Value = 2
$newSheet.Range('B4').Value = '160322'
$newSheet.Range('C4').Value = '华夏港股通精选股票（LOF）A'
$newSheet.Range('D4').Value = '13.82'
$newSheet.Range('E4').Value = '92.59'
$newSheet.Range('F4').Value = '2.40'
$newSheet.Range('G4').Value = '0.3317'
$newSheet.Range('H4').Value = 8
$newSheet.Range('A5').Value = 3
$newSheet.Range('B5').Value = '015663'
$newSheet.Range('C5').Value = '易米开鑫价值优选混合A'
$newSheet.Range('D5').Value = '2.97'
$newSheet.Range('E5').Value = '62.86'
$newSheet.Range('F5').Value = '2.92'
$newSheet.Range('G5').Value = '0.0867'
$newSheet.Range('H5').Value = 10
$newSheet.Range('A6').Value = 4
$newSheet.Range('B6').Value = '011114'
$newSheet.Range('C6').Value = '富国沪港深行业精选灵活配置混合C'
$newSheet.Range('D6').Value = '3.40'
$newSheet.Range('E6').Value = '79.74'
$newSheet.Range('F6').Value = '2.37'
$newSheet.Range('G6').Value = '0.0806'
$newSheet.Range('H6').Value = 8
$newSheet.Range('A7').Value = 5
$newSheet.Range('B7').Value = '006787'
$newSheet.Range('C7').Value = '泰康中证港股通大消费主题指数C'
$newSheet.Range('D7').Value = '1.28'
$newSheet.Range('E7').Value = '94.75'
$newSheet.Range('F7').Value = '3.34'
$newSheet.Range('G7').Value = '0.0428'
$newSheet.Range('H7').Value = 10
$newSheet.Range('A8').Value = 6
$newSheet.Range('B8').Value = '006786'
$newSheet.Range('C8').Value = '泰康中证港股通大消费主题指数A'
$newSheet.Range('D8').Value = '0.56'
$newSheet.Range('E8').Value = '94.75'
$newSheet.Range('F8').Value = '3.34'
$newSheet.Range('G8').Value = '0.0187'
$newSheet.Range('H8').Value = 10
$newSheet.Range('A9').Value = 7
$newSheet.Range('B9').Value = '012884'
$newSheet.Range('C9').Value = '华夏港股通精选股票（LOF）C'
$newSheet.Range('D9').Value = '0.69'
$newSheet.Range('E9').Value = '92.59'
$newSheet.Range('F9').Value = '2.40'
$newSheet.Range('G9').Value = '0.0166'
$newSheet.Range('H9').Value = 8
$newSheet.Range('A10').Value = 8
$newSheet.Range('B10').Value = '011534'
$newSheet.Range('C10').Value = '万家民瑞祥明6个月持有期混合型A'
$newSheet.Range('D10').Value = '1.83'
$newSheet.Range('E10').Value = '22.41'
$newSheet.Range('F10').Value = '0.87'
$newSheet.Range('G10').Value = '0.0159'
$newSheet.Range('H10').Value = 5
$newSheet.Range('A11').Value = 9
$newSheet.Range('B11').Value = '015664'
$newSheet.Range('C11').Value = '易米开鑫价值优选混合C'
$newSheet.Range('D11').Value = '0.50'
$newSheet.Range('E11').Value = '62.86'
$newSheet.Range('F11').Value = '2.92'
$newSheet.Range('G11').Value = '0.0146'
$newSheet.Range('H11').Value = 10
$newSheet.Range('A12').Value = 10
$newSheet.Range('B12').Value = '162416'
$newSheet.Range('C12').Value = '华宝港股通恒生香港35指数（LOF）'
$newSheet.Range('D12').Value = '0.25'
$newSheet.Range('E12').Value = '90.59'
$newSheet.Range('F12').Value = '3.47'
$newSheet.Range('G12').Value = '0.0087'
$newSheet.Range('H12').Value = 10
$newSheet.Range('A13').Value = 11
$newSheet.Range('B13').Value = '009734'
$newSheet.Range('C13').Value = '创金合信港股通大消费精选股票C'
$newSheet.Range('D13').Value = '0.19'
$newSheet.Range('E13').Value = '81.61'
$newSheet.Range('F13').Value = '3.44'
$newSheet.Range('G13').Value = '0.0065'
$newSheet.Range('H13').Value = 10
$newSheet.Range('A14').Value = 12
$newSheet.Range('B14').Value = '009733'
$newSheet.Range('C14').Value = '创金合信港股通大消费精选股票A'
$newSheet.Range('D14').Value = '0.09'
$newSheet.Range('E14').Value = '81.61'
$newSheet.Range('F14').Value = '3.44'
$newSheet.Range('G14').Value = '0.0031'
$newSheet.Range('H14').Value = 10
$newSheet.Range('A15').Value = 13
$newSheet.Range('B15').Value = '011535'
$newSheet.Range('C15').Value = '万家民瑞祥明6个月持有期混合型C'
$newSheet.Range('D15').Value = '0.19'
$newSheet.Range('E15').Value = '22.41'
$newSheet.Range('F15').Value = '0.87'
$newSheet.Range('G15').Value = '0.0017'
$newSheet.Range('H15').Value = 5

# ---- 2. Insert the 2022-Q4 row into "总计", shifting other rows down ----
$summary.Rows("2:2").Insert()

# Re-apply the index-column + plain-data style from the row immediately
# below (the former row 2, now row 3) so the new row matches the rest.
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)
$summary.Range("B3:D3").Copy()
$summary.Range("B2").PasteSpecial(-4122)

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 14
$summary.Range("D2").Value = 1.72

Write-Output "2022-Q4 sheet + summary row added"
